# Scheduled-runner update: refresh cached market-board price snapshots
# (currentAveragePrice / NQ / HQ) and the derived Leve profit columns
# across all eight crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 16570.666
$ws.Range("J3").Value = 16570.666
$ws.Range("L3").Value = 16570.666
$ws.Range("N3").Value = -16798.666

$ws.Range("H43").Value = 6291.6665
$ws.Range("I43").Value = 6312.5
$ws.Range("K43").Value = 6312.5
$ws.Range("M43").Value = -6243.5

$ws.Range("H61").Value = 550
$ws.Range("I61").Value = 550
$ws.Range("K61").Value = 1650
$ws.Range("M61").Value = -1478

$ws.Range("H76").Value = 3950
$ws.Range("I76").Value = 5500
$ws.Range("K76").Value = 5500
$ws.Range("M76").Value = -5185

$ws.Range("H79").Value = 3950
$ws.Range("I79").Value = 5500
$ws.Range("K79").Value = 5500
$ws.Range("M79").Value = -4408

$ws.Range("H87").Value = 74283.64
$ws.Range("J87").Value = 74283.64
$ws.Range("L87").Value = 74283.64
$ws.Range("N87").Value = -76779.64

$ws.Range("H90").Value = 74283.64
$ws.Range("J90").Value = 74283.64
$ws.Range("L90").Value = 222850.92
$ws.Range("N90").Value = -235330.92

$ws.Range("H93").Value = 313750
$ws.Range("J93").Value = 408333.34
$ws.Range("L93").Value = 408333.34
$ws.Range("N93").Value = -413325.34

$ws.Range("H94").Value = 503.25
$ws.Range("I94").Value = 503.25
$ws.Range("K94").Value = 503.25
$ws.Range("M94").Value = -52.25

$ws.Range("H102").Value = 16570.666
$ws.Range("J102").Value = 16570.666
$ws.Range("L102").Value = 16570.666
$ws.Range("N102").Value = -23060.666

$ws.Range("H132").Value = 3446.3914
$ws.Range("I132").Value = 3584.275
$ws.Range("K132").Value = 10752.825
$ws.Range("M132").Value = -8222.825000000001

$ws.Range("H135").Value = 1958
$ws.Range("I135").Value = 1737
$ws.Range("J135").Value = 2400
$ws.Range("K135").Value = 15633
$ws.Range("L135").Value = 21600
$ws.Range("M135").Value = -13098
$ws.Range("N135").Value = -26670

$ws.Range("H138").Value = 2264.5
$ws.Range("I138").Value = 997.5
$ws.Range("J138").Value = 2898
$ws.Range("K138").Value = 2992.5
$ws.Range("L138").Value = 8694
$ws.Range("M138").Value = 2147.5
$ws.Range("N138").Value = -18974

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2661.6
$ws.Range("I61").Value = 2661.6
$ws.Range("K61").Value = 2661.6
$ws.Range("M61").Value = -2449.6

$ws.Range("H88").Value = 2857.3333
$ws.Range("J88").Value = 2857.3333
$ws.Range("L88").Value = 2857.3333
$ws.Range("N88").Value = -3669.3333

$ws.Range("H91").Value = 2857.3333
$ws.Range("J91").Value = 2857.3333
$ws.Range("L91").Value = 2857.3333
$ws.Range("N91").Value = -5665.3333

$ws.Range("H92").Value = 150000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 150000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 150000
$ws.Range("M92").ClearContents() | Out-Null
$ws.Range("N92").Value = -154992

$ws.Range("H132").Value = 1465.1666
$ws.Range("I132").Value = 1411
$ws.Range("K132").Value = 4233
$ws.Range("M132").Value = -1703

$ws.Range("H136").Value = 2661.6
$ws.Range("I136").Value = 2661.6
$ws.Range("K136").Value = 7984.799999999999
$ws.Range("M136").Value = -5434.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 499
$ws.Range("I22").Value = 499
$ws.Range("K22").Value = 499
$ws.Range("M22").Value = -326

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents() | Out-Null

$ws.Range("H92").Value = 50000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 50000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 50000
$ws.Range("M92").ClearContents() | Out-Null
$ws.Range("N92").Value = -54992

$ws.Range("H134").Value = 2336.2666
$ws.Range("I134").Value = 2346.0715
$ws.Range("K134").Value = 7038.2145
$ws.Range("M134").Value = -4503.2145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 24000
$ws.Range("J74").Value = 22000
$ws.Range("L74").Value = 22000
$ws.Range("N74").Value = -23748

$ws.Range("H77").Value = 24000
$ws.Range("J77").Value = 22000
$ws.Range("L77").Value = 66000
$ws.Range("N77").Value = -74736

$ws.Range("H86").Value = 4299.5
$ws.Range("I86").Value = 4270.857
$ws.Range("K86").Value = 4270.857
$ws.Range("M86").Value = -3147.857

$ws.Range("H89").Value = 4299.5
$ws.Range("I89").Value = 4270.857
$ws.Range("K89").Value = 21354.285
$ws.Range("M89").Value = -15738.285

$ws.Range("H122").Value = 1732.7
$ws.Range("I122").Value = 1732.7
$ws.Range("K122").Value = 5198.1
$ws.Range("M122").Value = -2748.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1774008.6
$ws.Range("I4").Value = 1011144
$ws.Range("K4").Value = 3033432
$ws.Range("M4").Value = -3033320

$ws.Range("H51").Value = 1596
$ws.Range("J51").Value = 2888
$ws.Range("L51").Value = 8664
$ws.Range("N51").Value = -9584

$ws.Range("H54").Value = 4999
$ws.Range("J54").Value = 4999
$ws.Range("L54").Value = 14997
$ws.Range("N54").Value = -16115

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents() | Out-Null

$ws.Range("H58").Value = 4999
$ws.Range("J58").Value = 4999
$ws.Range("L58").Value = 14997
$ws.Range("N58").Value = -15253

$ws.Range("H137").Value = 4964.4
$ws.Range("J137").Value = 6833.3335
$ws.Range("L137").Value = 20500.0005
$ws.Range("N137").Value = -30700.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 240.85715
$ws.Range("I2").Value = 240.85715
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 240.85715
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -127.85715
$ws.Range("N2").ClearContents() | Out-Null

$ws.Range("H122").Value = 2970.0833
$ws.Range("I122").Value = 1612
$ws.Range("J122").Value = 4871.4
$ws.Range("K122").Value = 4836
$ws.Range("L122").Value = 14614.2
$ws.Range("M122").Value = -2386
$ws.Range("N122").Value = -19514.2

$ws.Range("H132").Value = 2751.4167
$ws.Range("I132").Value = 2461.8
$ws.Range("J132").Value = 4199.5
$ws.Range("K132").Value = 7385.400000000001
$ws.Range("L132").Value = 12598.5
$ws.Range("M132").Value = -4855.400000000001
$ws.Range("N132").Value = -17658.5

$ws.Range("H135").Value = 78000
$ws.Range("J135").Value = 78000
$ws.Range("L135").Value = 78000
$ws.Range("N135").Value = -88140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1508.5555
$ws.Range("I16").Value = 1096.5
$ws.Range("J16").Value = 2332.6667
$ws.Range("K16").Value = 1096.5
$ws.Range("L16").Value = 2332.6667
$ws.Range("M16").Value = -926.5
$ws.Range("N16").Value = -2672.6667

$ws.Range("H22").Value = 3107.6667
$ws.Range("I22").Value = 4632.6665
$ws.Range("K22").Value = 4632.6665
$ws.Range("M22").Value = -4337.6665

$ws.Range("H27").Value = 3107.6667
$ws.Range("I27").Value = 4632.6665
$ws.Range("K27").Value = 4632.6665
$ws.Range("M27").Value = -4525.6665

$ws.Range("H46").Value = 3335.0908
$ws.Range("I46").Value = 2424.75
$ws.Range("K46").Value = 2424.75
$ws.Range("M46").Value = -2236.75

$ws.Range("H68").Value = 3789
$ws.Range("I68").Value = 3937.1667
$ws.Range("K68").Value = 3937.1667
$ws.Range("M68").Value = -3188.1667

$ws.Range("H71").Value = 3789
$ws.Range("I71").Value = 3937.1667
$ws.Range("K71").Value = 19685.8335
$ws.Range("M71").Value = -15941.8335

$ws.Range("H128").Value = 82968.39999999999
$ws.Range("J128").Value = 82968.39999999999
$ws.Range("L128").Value = 82968.39999999999
$ws.Range("N128").Value = -92928.39999999999

$ws.Range("H132").Value = 2606.8
$ws.Range("I132").Value = 2786
$ws.Range("J132").Value = 2450
$ws.Range("K132").Value = 8358
$ws.Range("L132").Value = 7350
$ws.Range("M132").Value = -5828
$ws.Range("N132").Value = -12410

$ws.Range("H136").Value = 26319.8
$ws.Range("I136").Value = 13774.875
$ws.Range("J136").Value = 76499.5
$ws.Range("K136").Value = 41324.625
$ws.Range("L136").Value = 229498.5
$ws.Range("M136").Value = -38774.625
$ws.Range("N136").Value = -234598.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 10249.5
$ws.Range("J101").Value = 10249.5
$ws.Range("L101").Value = 10249.5
$ws.Range("N101").Value = -16739.5

$ws.Range("H126").Value = 2941
$ws.Range("I126").Value = 2897.5
$ws.Range("K126").Value = 8692.5
$ws.Range("M126").Value = -6222.5

$ws.Range("H132").Value = 3974.889
$ws.Range("I132").Value = 3974.889
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11924.667
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9394.667000000001
$ws.Range("N132").ClearContents() | Out-Null
